$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: Title (content unchanged, now bold) ---
$ws.Range("A1").Font.Bold = $true

# --- Row 2: header row gets a new (gold/Accent4) fill; B2 text changes ---
$ws.Range("B2").Value = "Properties to be Evaluated"
$ws.Range("A2:B2").Interior.ThemeColor = 8

# --- Row 3: pins-per-port sentence extended; new fill applied to body rows ---
$ws.Range("B3").Value = "To set the number of pins per port and also configure number of available ports."

# --- Row 4: direction register sentence, with "Direction Register" bold ---
$ws.Range("B4").Value = "To set the Direction Register as it can be configured as both input and output."
$ws.Range("B4").Characters(12, 18).Font.Bold = $true

# --- Row 5: function-select sentence, now with trailing period ---
$ws.Range("B5").Value = "To select the function of the port pin and set it for the general purpose or interrupt handling."

# --- Row 6: set Interrupt Configuration Register, bold + period ---
$ws.Range("B6").Value = "To set the Interrupt Configuration Register."
$ws.Range("B6").Characters(12, 34).Font.Bold = $true

# --- Row 7: check Interrupt Configuration Register, bold + period ---
$ws.Range("B7").Value = "To Check the Interrupt Configuration Register."
$ws.Range("B7").Characters(14, 34).Font.Bold = $true

# --- Apply the new body fill to the whole table (rows 2-7, both columns) ---
$ws.Range("A3:B7").Interior.ThemeColor = 8

# --- Selection moves to B2 ---
$ws.Range("B2").Select()
